$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

# Row 6 - signup test: lifeisgood user
$ws.Range("C6").Value = "lifeisgood@gmail.com"
$ws.Range("D6").Value = "hopscotch"
$ws.Range("E6").Value = "userManagement.Student"

# Row 7 - signup test: good user
$ws.Range("C7").Value = "good@gmail.com"
$ws.Range("D7").Value = "hopscotch"
$ws.Range("E7").Value = "Student"

# Row 8 - signup test: good user (duplicate)
$ws.Range("C8").Value = "good@gmail.com"
$ws.Range("D8").Value = "hopscotch"
$ws.Range("E8").Value = "Student"

# Row 9 - signup test: good user (duplicate)
$ws.Range("C9").Value = "good@gmail.com"
$ws.Range("D9").Value = "hopscotch"
$ws.Range("E9").Value = "Student"

# Row 10 - signup test: good user (duplicate)
$ws.Range("C10").Value = "good@gmail.com"
$ws.Range("D10").Value = "hopscotch"
$ws.Range("E10").Value = "Student"

# Row 11 - signup test: notgood user
$ws.Range("C11").Value = "notgood@gmail.com"
$ws.Range("D11").Value = "hopscotch"
$ws.Range("E11").Value = "Student"

# Row 12 - login test: hello user
$ws.Range("C12").Value = "hello@gmail.com"
$ws.Range("D12").Value = "hello"
$ws.Range("E12").Value = "Student"
$ws.Range("F12").Value = $true
